$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2 (ParticipantsTab) --------------------------------------------
# B2 (TabQuery) and C2 (StatQuery) keep their existing text; only D2/E2
# (TsvExcel / WebExcel file-name columns) change.
$ws.Range("D2").Value = "TC04_CDS_phs001524_SampleTumorStatus_NSD_TSVData.xlsx"
$ws.Range("E2").Value = "TC04_CDS_phs001524_SampleTumorStatus_NSD_WebData.xlsx"

# --- Row 3 (SamplesTab) --------------------------------------------------
# B3 query no longer selects the Tumor / Analyte Type columns.
$b3 = @"
SELECT
    DISTINCT (smp.sample_id) AS "Sample ID",
    sp.participant_id AS "Participant ID", 
    s.study_name AS "Study Name",
    s.phs_accession AS Accession
FROM 
    df_participant sp
JOIN 
    df_study s ON sp."study.phs_accession" = s.phs_accession
JOIN 
    df_sample smp ON smp."participant.study_participant_id" = sp.study_participant_id
JOIN
    df_diagnosis d ON d."participant.study_participant_id" = sp.study_participant_id
JOIN
    df_program p ON p.program_acronym = s."program.program_acronym"
JOIN
    df_file f1 ON f1."sample.sample_id" = smp.sample_id
JOIN
    df_genomic_info gi ON gi."file.file_id" = f1.file_id
WHERE 
    s.phs_accession = 'phs001524' AND smp.sample_tumor_status = 'Not specified in data'
ORDER BY 
    smp.sample_id ASC
LIMIT 100;
"@
$ws.Range("B3").Value = $b3

# The TsvExcel/WebExcel file-name columns no longer apply to this row.
$ws.Range("D3").ClearContents()
$ws.Range("E3").ClearContents()

# --- Row 4 (FilesTab) -----------------------------------------------------
# B4 query text is unchanged; its TsvExcel/WebExcel columns are removed.
$ws.Range("D4").ClearContents()
$ws.Range("E4").ClearContents()

# --- Selection / scroll position ------------------------------------------
$win = $excel.ActiveWindow
$win.ScrollRow = 3
$win.ScrollColumn = 1
$ws.Range("C3").Select() | Out-Null
